$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 26.56908266666666
$ws.Cells.Item(2, 8).Value = 79.70724799999999
$ws.Cells.Item(2, 9).Value = 0.7506383589294218
$ws.Cells.Item(2, 10).Value = 0.7506383589294219
$ws.Cells.Item(2, 13).Value = 2.231113333333334
$ws.Cells.Item(2, 14).Value = 6.69334
$ws.Cells.Item(2, 15).Value = 0.01598125358798882
$ws.Cells.Item(2, 16).Value = 0.01598125358798882
$ws.Cells.Item(2, 17).Value = 59.27863459203555
$ws.Cells.Item(2, 18).Value = 533.5077113283199
$ws.Cells.Item(2, 19).Value = 0.01199614196692286
$ws.Cells.Item(2, 20).Value = 0.01199614196692286
$ws.Cells.Item(3, 7).Value = 26.56908266666666
$ws.Cells.Item(3, 8).Value = 79.70724799999999
$ws.Cells.Item(3, 9).Value = 0.7506383589294218
$ws.Cells.Item(3, 10).Value = 0.7506383589294219
$ws.Cells.Item(3, 15).Value = 0.1634493267640196
$ws.Cells.Item(3, 16).Value = 0.1634493267640195
$ws.Cells.Item(3, 17).Value = 606.2761511300106
$ws.Cells.Item(3, 18).Value = 5456.485360170095
$ws.Cells.Item(3, 19).Value = 0.1226913344102625
$ws.Cells.Item(3, 20).Value = 0.1226913344102625
$ws.Cells.Item(4, 7).Value = 26.56908266666666
$ws.Cells.Item(4, 8).Value = 79.70724799999999
$ws.Cells.Item(4, 9).Value = 0.7506383589294218
$ws.Cells.Item(4, 10).Value = 0.7506383589294219
$ws.Cells.Item(4, 13).Value = 58.02175166666666
$ws.Cells.Item(4, 14).Value = 174.065255
$ws.Cells.Item(4, 15).Value = 0.4156043142904646
$ws.Cells.Item(4, 16).Value = 0.4156043142904646
$ws.Cells.Item(4, 17).Value = 1541.584716496471
$ws.Cells.Item(4, 18).Value = 13874.26244846824
$ws.Cells.Item(4, 19).Value = 0.3119685404429819
$ws.Cells.Item(4, 20).Value = 0.311968540442982
$ws.Cells.Item(5, 7).Value = 26.56908266666666
$ws.Cells.Item(5, 8).Value = 79.70724799999999
$ws.Cells.Item(5, 9).Value = 0.7506383589294218
$ws.Cells.Item(5, 10).Value = 0.7506383589294219
$ws.Cells.Item(5, 13).Value = 15.16934033333333
$ws.Cells.Item(5, 14).Value = 45.508021
$ws.Cells.Item(5, 15).Value = 0.1086565487318021
$ws.Cells.Item(5, 16).Value = 0.1086565487318021
$ws.Cells.Item(5, 17).Value = 403.0354573151341
$ws.Cells.Item(5, 18).Value = 3627.319115836207
$ws.Cells.Item(5, 19).Value = 0.0815617734269747
$ws.Cells.Item(5, 20).Value = 0.0815617734269747
$ws.Cells.Item(6, 7).Value = 26.56908266666666
$ws.Cells.Item(6, 8).Value = 79.70724799999999
$ws.Cells.Item(6, 9).Value = 0.7506383589294218
$ws.Cells.Item(6, 10).Value = 0.7506383589294219
$ws.Cells.Item(6, 13).Value = 41.36709099999999
$ws.Cells.Item(6, 14).Value = 124.101273
$ws.Cells.Item(6, 15).Value = 0.2963085566257249
$ws.Cells.Item(6, 16).Value = 0.2963085566257249
$ws.Cells.Item(6, 17).Value = 1099.085660458522
$ws.Cells.Item(6, 18).Value = 9891.770944126702
$ws.Cells.Item(6, 19).Value = 0.2224205686822798
$ws.Cells.Item(6, 20).Value = 0.2224205686822798
$ws.Cells.Item(7, 9).Value = 0.004944072121179194
$ws.Cells.Item(7, 10).Value = 0.004944072121179195
$ws.Cells.Item(7, 13).Value = 2.231113333333334
$ws.Cells.Item(7, 14).Value = 6.69334
$ws.Cells.Item(7, 15).Value = 0.01598125358798882
$ws.Cells.Item(7, 16).Value = 0.01598125358798882
$ws.Cells.Item(7, 17).Value = 0.3904381399933333
$ws.Cells.Item(7, 18).Value = 3.51394325994
$ws.Cells.Item(7, 19).Value = 0.00007901247032587049
$ws.Cells.Item(7, 20).Value = 0.0000790124703258705
$ws.Cells.Item(8, 9).Value = 0.004944072121179194
$ws.Cells.Item(8, 10).Value = 0.004944072121179195
$ws.Cells.Item(8, 15).Value = 0.1634493267640196
$ws.Cells.Item(8, 16).Value = 0.1634493267640195
$ws.Cells.Item(8, 19).Value = 0.0008081052596794975
$ws.Cells.Item(8, 20).Value = 0.0008081052596794975
$ws.Cells.Item(9, 9).Value = 0.004944072121179194
$ws.Cells.Item(9, 10).Value = 0.004944072121179195
$ws.Cells.Item(9, 13).Value = 58.02175166666666
$ws.Cells.Item(9, 14).Value = 174.065255
$ws.Cells.Item(9, 15).Value = 0.4156043142904646
$ws.Cells.Item(9, 16).Value = 0.4156043142904646
$ws.Cells.Item(9, 17).Value = 10.15363247641166
$ws.Cells.Item(9, 18).Value = 91.38269228770498
$ws.Cells.Item(9, 19).Value = 0.002054777703725282
$ws.Cells.Item(9, 20).Value = 0.002054777703725282
$ws.Cells.Item(10, 9).Value = 0.004944072121179194
$ws.Cells.Item(10, 10).Value = 0.004944072121179195
$ws.Cells.Item(10, 13).Value = 15.16934033333333
$ws.Cells.Item(10, 14).Value = 45.508021
$ws.Cells.Item(10, 15).Value = 0.1086565487318021
$ws.Cells.Item(10, 16).Value = 0.1086565487318021
$ws.Cells.Item(10, 17).Value = 2.654589050312333
$ws.Cells.Item(10, 18).Value = 23.891301452811
$ws.Cells.Item(10, 19).Value = 0.0005372058133684514
$ws.Cells.Item(10, 20).Value = 0.0005372058133684515
$ws.Cells.Item(11, 9).Value = 0.004944072121179194
$ws.Cells.Item(11, 10).Value = 0.004944072121179195
$ws.Cells.Item(11, 13).Value = 41.36709099999999
$ws.Cells.Item(11, 14).Value = 124.101273
$ws.Cells.Item(11, 15).Value = 0.2963085566257249
$ws.Cells.Item(11, 16).Value = 0.2963085566257249
$ws.Cells.Item(11, 17).Value = 7.239116823726999
$ws.Cells.Item(11, 18).Value = 65.15205141354299
$ws.Cells.Item(11, 19).Value = 0.001464970874080093
$ws.Cells.Item(11, 20).Value = 0.001464970874080093
$ws.Cells.Item(12, 7).Value = 5.094400666666666
$ws.Cells.Item(12, 8).Value = 15.283202
$ws.Cells.Item(12, 9).Value = 0.1439286634067062
$ws.Cells.Item(12, 10).Value = 0.1439286634067062
$ws.Cells.Item(12, 13).Value = 2.231113333333334
$ws.Cells.Item(12, 14).Value = 6.69334
$ws.Cells.Item(12, 15).Value = 0.01598125358798882
$ws.Cells.Item(12, 16).Value = 0.01598125358798882
$ws.Cells.Item(12, 17).Value = 11.36618525274222
$ws.Cells.Item(12, 18).Value = 102.29566727468
$ws.Cells.Item(12, 19).Value = 0.002300160468482859
$ws.Cells.Item(12, 20).Value = 0.002300160468482859
$ws.Cells.Item(13, 7).Value = 5.094400666666666
$ws.Cells.Item(13, 8).Value = 15.283202
$ws.Cells.Item(13, 9).Value = 0.1439286634067062
$ws.Cells.Item(13, 10).Value = 0.1439286634067062
$ws.Cells.Item(13, 15).Value = 0.1634493267640196
$ws.Cells.Item(13, 16).Value = 0.1634493267640195
$ws.Cells.Item(13, 17).Value = 116.2484105021726
$ws.Cells.Item(13, 18).Value = 1046.235694519554
$ws.Cells.Item(13, 19).Value = 0.02352504313587131
$ws.Cells.Item(13, 20).Value = 0.02352504313587131
$ws.Cells.Item(14, 7).Value = 5.094400666666666
$ws.Cells.Item(14, 8).Value = 15.283202
$ws.Cells.Item(14, 9).Value = 0.1439286634067062
$ws.Cells.Item(14, 10).Value = 0.1439286634067062
$ws.Cells.Item(14, 13).Value = 58.02175166666666
$ws.Cells.Item(14, 14).Value = 174.065255
$ws.Cells.Item(14, 15).Value = 0.4156043142904646
$ws.Cells.Item(14, 16).Value = 0.4156043142904646
$ws.Cells.Item(14, 17).Value = 295.5860503718344
$ws.Cells.Item(14, 18).Value = 2660.27445334651
$ws.Cells.Item(14, 19).Value = 0.05981737346188722
$ws.Cells.Item(14, 20).Value = 0.05981737346188722
$ws.Cells.Item(15, 7).Value = 5.094400666666666
$ws.Cells.Item(15, 8).Value = 15.283202
$ws.Cells.Item(15, 9).Value = 0.1439286634067062
$ws.Cells.Item(15, 10).Value = 0.1439286634067062
$ws.Cells.Item(15, 13).Value = 15.16934033333333
$ws.Cells.Item(15, 14).Value = 45.508021
$ws.Cells.Item(15, 15).Value = 0.1086565487318021
$ws.Cells.Item(15, 16).Value = 0.1086565487318021
$ws.Cells.Item(15, 17).Value = 77.27869750702688
$ws.Cells.Item(15, 18).Value = 695.5082775632419
$ws.Cells.Item(15, 19).Value = 0.01563879182935392
$ws.Cells.Item(15, 20).Value = 0.01563879182935392
$ws.Cells.Item(16, 7).Value = 5.094400666666666
$ws.Cells.Item(16, 8).Value = 15.283202
$ws.Cells.Item(16, 9).Value = 0.1439286634067062
$ws.Cells.Item(16, 10).Value = 0.1439286634067062
$ws.Cells.Item(16, 13).Value = 41.36709099999999
$ws.Cells.Item(16, 14).Value = 124.101273
$ws.Cells.Item(16, 15).Value = 0.2963085566257249
$ws.Cells.Item(16, 16).Value = 0.2963085566257249
$ws.Cells.Item(16, 17).Value = 210.7405359684606
$ws.Cells.Item(16, 18).Value = 1896.664823716146
$ws.Cells.Item(16, 19).Value = 0.04264729451111091
$ws.Cells.Item(16, 20).Value = 0.04264729451111092
$ws.Cells.Item(17, 7).Value = 0.05498833333333333
$ws.Cells.Item(17, 8).Value = 0.164965
$ws.Cells.Item(17, 9).Value = 0.001553548265532792
$ws.Cells.Item(17, 10).Value = 0.001553548265532792
$ws.Cells.Item(17, 13).Value = 2.231113333333334
$ws.Cells.Item(17, 14).Value = 6.69334
$ws.Cells.Item(17, 15).Value = 0.01598125358798882
$ws.Cells.Item(17, 16).Value = 0.01598125358798882
$ws.Cells.Item(17, 17).Value = 0.1226852036777778
$ws.Cells.Item(17, 18).Value = 1.1041668331
$ws.Cells.Item(17, 19).Value = 0.00002482764879265973
$ws.Cells.Item(17, 20).Value = 0.00002482764879265973
$ws.Cells.Item(18, 7).Value = 0.05498833333333333
$ws.Cells.Item(18, 8).Value = 0.164965
$ws.Cells.Item(18, 9).Value = 0.001553548265532792
$ws.Cells.Item(18, 10).Value = 0.001553548265532792
$ws.Cells.Item(18, 15).Value = 0.1634493267640196
$ws.Cells.Item(18, 16).Value = 0.1634493267640195
$ws.Cells.Item(18, 17).Value = 1.254771024978333
$ws.Cells.Item(18, 18).Value = 11.292939224805
$ws.Cells.Item(18, 19).Value = 0.0002539264180967451
$ws.Cells.Item(18, 20).Value = 0.0002539264180967451
$ws.Cells.Item(19, 7).Value = 0.05498833333333333
$ws.Cells.Item(19, 8).Value = 0.164965
$ws.Cells.Item(19, 9).Value = 0.001553548265532792
$ws.Cells.Item(19, 10).Value = 0.001553548265532792
$ws.Cells.Item(19, 13).Value = 58.02175166666666
$ws.Cells.Item(19, 14).Value = 174.065255
$ws.Cells.Item(19, 15).Value = 0.4156043142904646
$ws.Cells.Item(19, 16).Value = 0.4156043142904646
$ws.Cells.Item(19, 17).Value = 3.190519421230555
$ws.Cells.Item(19, 18).Value = 28.714674791075
$ws.Cells.Item(19, 19).Value = 0.0006456613616138964
$ws.Cells.Item(19, 20).Value = 0.0006456613616138966
$ws.Cells.Item(20, 7).Value = 0.05498833333333333
$ws.Cells.Item(20, 8).Value = 0.164965
$ws.Cells.Item(20, 9).Value = 0.001553548265532792
$ws.Cells.Item(20, 10).Value = 0.001553548265532792
$ws.Cells.Item(20, 13).Value = 15.16934033333333
$ws.Cells.Item(20, 14).Value = 45.508021
$ws.Cells.Item(20, 15).Value = 0.1086565487318021
$ws.Cells.Item(20, 16).Value = 0.1086565487318021
$ws.Cells.Item(20, 17).Value = 0.834136742696111
$ws.Cells.Item(20, 18).Value = 7.507230684265
$ws.Cells.Item(20, 19).Value = 0.0001688031928210704
$ws.Cells.Item(20, 20).Value = 0.0001688031928210704
$ws.Cells.Item(21, 7).Value = 0.05498833333333333
$ws.Cells.Item(21, 8).Value = 0.164965
$ws.Cells.Item(21, 9).Value = 0.001553548265532792
$ws.Cells.Item(21, 10).Value = 0.001553548265532792
$ws.Cells.Item(21, 13).Value = 41.36709099999999
$ws.Cells.Item(21, 14).Value = 124.101273
$ws.Cells.Item(21, 15).Value = 0.2963085566257249
$ws.Cells.Item(21, 16).Value = 0.2963085566257249
$ws.Cells.Item(21, 17).Value = 2.274707388938333
$ws.Cells.Item(21, 18).Value = 20.472366500445
$ws.Cells.Item(21, 19).Value = 0.0004603296442084199
$ws.Cells.Item(21, 20).Value = 0.0004603296442084199
$ws.Cells.Item(22, 7).Value = 3.501848333333333
$ws.Cells.Item(22, 8).Value = 10.505545
$ws.Cells.Item(22, 9).Value = 0.09893535727715995
$ws.Cells.Item(22, 10).Value = 0.09893535727715996
$ws.Cells.Item(22, 13).Value = 2.231113333333334
$ws.Cells.Item(22, 14).Value = 6.69334
$ws.Cells.Item(22, 15).Value = 0.01598125358798882
$ws.Cells.Item(22, 16).Value = 0.01598125358798882
$ws.Cells.Item(22, 17).Value = 7.813020507811111
$ws.Cells.Item(22, 18).Value = 70.3171845703
$ws.Cells.Item(22, 19).Value = 0.001581111033464568
$ws.Cells.Item(22, 20).Value = 0.001581111033464568
$ws.Cells.Item(23, 7).Value = 3.501848333333333
$ws.Cells.Item(23, 8).Value = 10.505545
$ws.Cells.Item(23, 9).Value = 0.09893535727715995
$ws.Cells.Item(23, 10).Value = 0.09893535727715996
$ws.Cells.Item(23, 15).Value = 0.1634493267640196
$ws.Cells.Item(23, 16).Value = 0.1634493267640195
$ws.Cells.Item(23, 17).Value = 79.90818335771833
$ws.Cells.Item(23, 18).Value = 719.1736502194649
$ws.Cells.Item(23, 19).Value = 0.01617091754010954
$ws.Cells.Item(23, 20).Value = 0.01617091754010954
$ws.Cells.Item(24, 7).Value = 3.501848333333333
$ws.Cells.Item(24, 8).Value = 10.505545
$ws.Cells.Item(24, 9).Value = 0.09893535727715995
$ws.Cells.Item(24, 10).Value = 0.09893535727715996
$ws.Cells.Item(24, 13).Value = 58.02175166666666
$ws.Cells.Item(24, 14).Value = 174.065255
$ws.Cells.Item(24, 15).Value = 0.4156043142904646
$ws.Cells.Item(24, 16).Value = 0.4156043142904646
$ws.Cells.Item(24, 17).Value = 203.1833743709972
$ws.Cells.Item(24, 18).Value = 1828.650369338975
$ws.Cells.Item(24, 19).Value = 0.04111796132025618
$ws.Cells.Item(24, 20).Value = 0.04111796132025619
$ws.Cells.Item(25, 7).Value = 3.501848333333333
$ws.Cells.Item(25, 8).Value = 10.505545
$ws.Cells.Item(25, 9).Value = 0.09893535727715995
$ws.Cells.Item(25, 10).Value = 0.09893535727715996
$ws.Cells.Item(25, 13).Value = 15.16934033333333
$ws.Cells.Item(25, 14).Value = 45.508021
$ws.Cells.Item(25, 15).Value = 0.1086565487318021
$ws.Cells.Item(25, 16).Value = 0.1086565487318021
$ws.Cells.Item(25, 17).Value = 53.12072916404944
$ws.Cells.Item(25, 18).Value = 478.086562476445
$ws.Cells.Item(25, 19).Value = 0.01074997446928398
$ws.Cells.Item(25, 20).Value = 0.01074997446928399
$ws.Cells.Item(26, 7).Value = 3.501848333333333
$ws.Cells.Item(26, 8).Value = 10.505545
$ws.Cells.Item(26, 9).Value = 0.09893535727715995
$ws.Cells.Item(26, 10).Value = 0.09893535727715996
$ws.Cells.Item(26, 13).Value = 41.36709099999999
$ws.Cells.Item(26, 14).Value = 124.101273
$ws.Cells.Item(26, 15).Value = 0.2963085566257249
$ws.Cells.Item(26, 16).Value = 0.2963085566257249
$ws.Cells.Item(26, 17).Value = 144.8612786731983
$ws.Cells.Item(26, 18).Value = 1303.751508058785
$ws.Cells.Item(26, 19).Value = 0.02931539291404568
$ws.Cells.Item(26, 20).Value = 0.02931539291404568
